$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Negative-word column (A): anchor-word list re-ranked after re-running
#     the toy-spam confidence report with a minimum-occurrence filter of 5 ---
$ws.Range("A6").Value = "disappointed"
$ws.Range("A7").Value = "poor"
$ws.Range("A8").Value = "broke"
$ws.Range("A9").Value = "waste"
$ws.Range("A10").Value = "smaller"
$ws.Range("A11").Value = "junk"
$ws.Range("A13").Value = "paint"
$ws.Range("A14").Value = "broken"
$ws.Range("A16").Value = "difficult"
$ws.Range("A17").Value = "apart"
$ws.Range("A20").Value = "thought"
$ws.Range("A21").Value = "size"
$ws.Range("A22").Value = "would"
$ws.Range("A23").Value = "item"
$ws.Range("A24").Value = "money"
$ws.Range("A25").Value = "work"
$ws.Range("A27").Value = "price"

# --- Positive-word column (J): same re-rank, only rows 3-4 swap order ---
$ws.Range("J3").Value = "awesome"
$ws.Range("J4").Value = "wonderful"

# --- Updated statistics for the negative-word table (columns B:H, rows 3-28) ---
$ws.Range("B3").Value = 0.9565217391304348
$ws.Range("C3").Value = 44
$ws.Range("D3").Value = 44
$ws.Range("H3").Value = 2
$ws.Range("B4").Value = 0.8863636363636364
$ws.Range("C4").Value = 39
$ws.Range("D4").Value = 39
$ws.Range("H4").Value = 5
$ws.Range("B5").Value = 0.765625
$ws.Range("C5").Value = 49
$ws.Range("D5").Value = 49
$ws.Range("H5").Value = 15
$ws.Range("B6").Value = 0.7311827956989247
$ws.Range("C6").Value = 136
$ws.Range("D6").Value = 136
$ws.Range("H6").Value = 50
$ws.Range("B7").Value = 0.704225352112676
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = 50
$ws.Range("H7").Value = 21
$ws.Range("B8").Value = 0.7038834951456311
$ws.Range("C8").Value = 145
$ws.Range("D8").Value = 145
$ws.Range("H8").Value = 61
$ws.Range("B9").Value = 0.6621621621621622
$ws.Range("C9").Value = 98
$ws.Range("D9").Value = 98
$ws.Range("H9").Value = 50
$ws.Range("B10").Value = 0.5966386554621849
$ws.Range("C10").Value = 71
$ws.Range("D10").Value = 71
$ws.Range("H10").Value = 48
$ws.Range("B11").Value = 0.5454545454545454
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 30
$ws.Range("H11").Value = 25
$ws.Range("B12").Value = 0.4869565217391305
$ws.Range("C12").Value = 168
$ws.Range("D12").Value = 168
$ws.Range("H12").Value = 177
$ws.Range("B13").Value = 0.4761904761904762
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = 30
$ws.Range("H13").Value = 33
$ws.Range("B14").Value = 0.4216867469879518
$ws.Range("C14").Value = 35
$ws.Range("D14").Value = 35
$ws.Range("H14").Value = 48
$ws.Range("B15").Value = 0.4173228346456693
$ws.Range("C15").Value = 53
$ws.Range("D15").Value = 53
$ws.Range("H15").Value = 74
$ws.Range("B16").Value = 0.3483146067415731
$ws.Range("C16").Value = 31
$ws.Range("D16").Value = 31
$ws.Range("H16").Value = 58
$ws.Range("B17").Value = 0.3473684210526316
$ws.Range("C17").Value = 33
$ws.Range("D17").Value = 33
$ws.Range("H17").Value = 62
$ws.Range("B18").Value = 0.3203125
$ws.Range("C18").Value = 41
$ws.Range("D18").Value = 41
$ws.Range("H18").Value = 87
$ws.Range("B19").Value = 0.2938388625592417
$ws.Range("C19").Value = 62
$ws.Range("D19").Value = 62
$ws.Range("H19").Value = 149
$ws.Range("B20").Value = 0.2920792079207921
$ws.Range("C20").Value = 59
$ws.Range("D20").Value = 59
$ws.Range("H20").Value = 143
$ws.Range("B21").Value = 0.2474226804123711
$ws.Range("C21").Value = 48
$ws.Range("D21").Value = 48
$ws.Range("H21").Value = 146
$ws.Range("B22").Value = 0.1899109792284867
$ws.Range("C22").Value = 128
$ws.Range("D22").Value = 128
$ws.Range("H22").Value = 546
$ws.Range("B23").Value = 0.1884057971014493
$ws.Range("C23").Value = 52
$ws.Range("D23").Value = 52
$ws.Range("H23").Value = 224
$ws.Range("B24").Value = 0.1708860759493671
$ws.Range("C24").Value = 54
$ws.Range("D24").Value = 54
$ws.Range("H24").Value = 262
$ws.Range("B25").Value = 0.1708860759493671
$ws.Range("C25").Value = 54
$ws.Range("D25").Value = 54
$ws.Range("H25").Value = 262
$ws.Range("B26").Value = 0.1409691629955947
$ws.Range("C26").Value = 64
$ws.Range("D26").Value = 64
$ws.Range("H26").Value = 390
$ws.Range("B27").Value = 0.138328530259366
$ws.Range("C27").Value = 48
$ws.Range("D27").Value = 49
$ws.Range("E27").Value = 0.02
$ws.Range("F27").Value = 0.98
$ws.Range("G27").Value = $true
$ws.Range("H27").Value = 299
$ws.Range("B28").Value = 0.05921052631578947
$ws.Range("C28").Value = 36
$ws.Range("D28").Value = 36
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = $false
$ws.Range("H28").Value = 572

# --- Updated statistics for the positive-word table (columns K:Q, rows 3-14) ---
$ws.Range("K3").Value = 0.8
$ws.Range("L3").Value = 52
$ws.Range("M3").Value = 52
$ws.Range("Q3").Value = 13
$ws.Range("K4").Value = 0.7857142857142857
$ws.Range("L4").Value = 44
$ws.Range("M4").Value = 44
$ws.Range("Q4").Value = 12
$ws.Range("K5").Value = 0.6129032258064516
$ws.Range("L5").Value = 57
$ws.Range("M5").Value = 57
$ws.Range("Q5").Value = 36
$ws.Range("K6").Value = 0.6037735849056604
$ws.Range("L6").Value = 32
$ws.Range("M6").Value = 32
$ws.Range("Q6").Value = 21
$ws.Range("K7").Value = 0.5
$ws.Range("L7").Value = 32
$ws.Range("M7").Value = 32
$ws.Range("Q7").Value = 32
$ws.Range("K8").Value = 0.3352459016393443
$ws.Range("L8").Value = 409
$ws.Range("M8").Value = 409
$ws.Range("Q8").Value = 811
$ws.Range("K9").Value = 0.3113342898134864
$ws.Range("L9").Value = 217
$ws.Range("M9").Value = 217
$ws.Range("Q9").Value = 480
$ws.Range("K10").Value = 0.2551867219917012
$ws.Range("L10").Value = 123
$ws.Range("M10").Value = 123
$ws.Range("Q10").Value = 359
$ws.Range("K11").Value = 0.1867469879518072
$ws.Range("L11").Value = 31
$ws.Range("M11").Value = 31
$ws.Range("Q11").Value = 135
$ws.Range("K12").Value = 0.1743119266055046
$ws.Range("L12").Value = 57
$ws.Range("M12").Value = 57
$ws.Range("Q12").Value = 270
$ws.Range("K13").Value = 0.08501314636283962
$ws.Range("L13").Value = 97
$ws.Range("M13").Value = 97
$ws.Range("Q13").Value = 1044
$ws.Range("K14").Value = 0.02790395846852693
$ws.Range("M14").Value = 43
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 1498
